# Apply the "(Why is it a bad model?)" title addition plus the two new
# blank paragraphs that follow it.
#
# The title paragraph currently holds a single run:
#   "Bitcoin stock-to-flow model"
# The target adds two more runs (kept distinct, matching the authored
# OOXML): a literal space, then "(Why is it a bad model?)" - both sz/szCs
# 48 (half-points), same as the existing run. We rebuild the paragraph's
# content via Range.InsertXML (WordprocessingML package payload) so the
# three runs are emitted as separate <w:r> elements instead of being
# coalesced into one.
#
# After that we insert two new empty paragraphs directly below the title:
#   1) blank line styled like the existing empty paragraph (sz/szCs 24)
#   2) centered blank paragraph styled sz/szCs 48 (no runs)

$d = $word.ActiveDocument

# Locate the title paragraph ("Bitcoin stock-to-flow model") rather than
# assuming a fixed index, in case the document contains leading content.
$titleIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Bitcoin stock-to-flow model*") {
        $titleIndex = $i
        break
    }
}

function New-WordPackageXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1: extend the title paragraph with the two extra runs -----------
$titlePara = $d.Paragraphs($titleIndex)
$titleRange = $titlePara.Range
$titleRange.End = $titleRange.End - 1   # exclude the paragraph mark

$titleBodyXml =
    '<w:p>' +
      '<w:pPr>' +
        '<w:jc w:val="center"/>' +
        '<w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr>' +
      '</w:pPr>' +
      '<w:r w:rsidRPr="00BB0D04">' +
        '<w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr>' +
        '<w:t>Bitcoin stock-to-flow model</w:t>' +
      '</w:r>' +
      '<w:r>' +
        '<w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr>' +
        '<w:t xml:space="preserve"> </w:t>' +
      '</w:r>' +
      '<w:r>' +
        '<w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr>' +
        '<w:t>(Why is it a bad model?)</w:t>' +
      '</w:r>' +
    '</w:p>'

$titleRange.InsertXML((New-WordPackageXml $titleBodyXml))

# --- Step 2: insert a blank paragraph (sz 24) right after the title -------
$titlePara = $d.Paragraphs($titleIndex)
$afterTitle = $titlePara.Range
$afterTitle.Collapse(0)
$afterTitle.InsertParagraphAfter()

$blankLine = $d.Paragraphs($titleIndex + 1).Range
$blankLineBodyXml =
    '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$blankLine.InsertXML((New-WordPackageXml $blankLineBodyXml))

# --- Step 3: insert a centered blank paragraph (sz 48) after that ---------
$blankPara = $d.Paragraphs($titleIndex + 1)
$afterBlank = $blankPara.Range
$afterBlank.Collapse(0)
$afterBlank.InsertParagraphAfter()

$centeredBlank = $d.Paragraphs($titleIndex + 2).Range
$centeredBlankBodyXml =
    '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr></w:p>'
$centeredBlank.InsertXML((New-WordPackageXml $centeredBlankBodyXml))
